$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the connection-string keys for clarity (values, not just shared-string text)
$ws.Range("A1").Value = "appDbConStr"
$ws.Range("B1").Value = "connStr"
$ws.Range("A2").Value = "reportsConStr"
$ws.Range("B2").Value = "connStr"

# Column A: auto-fit to the new (longer) key names
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666

# Column B: manually narrowed down from its old (legacy) width
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666

# Update the active selection
$ws.Range("E7").Select()
